$wb = $excel.ActiveWorkbook

# --- Update selection on the "Test" sheet ---
$wsTest = $wb.Worksheets.Item("Test")
$wsTest.Range("D5:G10").Select()

# --- Update selection on the "VariableNames" sheet (also drops tabSelected there) ---
$wsVar = $wb.Worksheets.Item("VariableNames")
$wsVar.Range("E16").Select()

# --- Add the new "Simplify" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Simplify"

# Header row + first vertical block
$ws.Range("C5").Value = "A"
$ws.Range("F5").Value = "A"
$ws.Range("G5").Value = "B"
$ws.Range("H5").Value = "C"
$ws.Range("I5").Value = "D"

$ws.Range("C6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 4

$ws.Range("C7").Value = 2
$ws.Range("C8").Value = 3
$ws.Range("C9").Value = 4

$ws.Range("C10").Value = 5
$ws.Range("F10").Value = "A"
$ws.Range("G10").Value = "B"
$ws.Range("H10").Value = "C"
$ws.Range("I10").Value = "D"

$ws.Range("C11").Value = 6
$ws.Range("F11").Value = $true
$ws.Range("G11").Value = $false
$ws.Range("H11").Value = $false
$ws.Range("I11").Value = $true

$ws.Range("C12").Value = 7
$ws.Range("C13").Value = 8
$ws.Range("C14").Value = 9

$ws.Range("C15").Value = 10
$ws.Range("F15").Value = "A"

$ws.Range("F16").Value = "one"
$ws.Range("F17").Value = "two"
$ws.Range("F18").Value = "three"
$ws.Range("F19").Value = "four"
$ws.Range("F20").Value = "five"

$ws.Range("G19").Select()

# --- Defined names for the new region blocks ---
$wb.Names.Add("Simplify1", '=Simplify!$C$5:$C$15')
$wb.Names.Add("Simplify2", '=Simplify!$F$5:$I$6')
$wb.Names.Add("Simplify3", '=Simplify!$F$10:$I$11')
$wb.Names.Add("Simplify4", '=Simplify!$F$15:$F$20')
